$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5

# Row 3
$ws.Range("G3").Value = 1.83
$ws.Range("I3").Value = 5.5
$ws.Range("J3").Value = 2.6
$ws.Range("K3").Value = 1.95
$ws.Range("AL3").Value = 67
$ws.Range("AO3").Value = 10
$ws.Range("BA3").Value = 201

# Row 4
$ws.Range("G4").Value = 3.4
$ws.Range("I4").Value = 2.38
$ws.Range("J4").Value = 4
$ws.Range("L4").Value = 3.25
$ws.Range("W4").Value = 7.5
$ws.Range("Z4").Value = 34
$ws.Range("AB4").Value = 41
$ws.Range("AH4").Value = 10
$ws.Range("AJ4").Value = 23
$ws.Range("AQ4").Value = 67
$ws.Range("AW4").Value = 4.33

# Row 5
$ws.Range("G5").Value = 1.55
$ws.Range("H5").Value = 3.9
$ws.Range("I5").Value = 6.25
$ws.Range("J5").Value = 2.2
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("Y5").Value = 9
$ws.Range("AF5").Value = 81
$ws.Range("AS5").Value = 201

# Row 7
$ws.Range("G7").Value = 2.37
$ws.Range("H7").Value = 3.1
$ws.Range("M7").Value = 1.08
$ws.Range("X7").Value = 10.75
$ws.Range("Z7").Value = 24
$ws.Range("AJ7").Value = 40
